$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Operation Test")

# ---------------------------------------------------------------
# 1. Duplicate an existing test-case group (rows 26:32, the "2."
#    group) down onto the currently-blank rows 70:76 to create the
#    new "4. Dang ky tai khoan" group, keeping the exact same
#    layout/merges/styles as the other groups on the sheet.
# ---------------------------------------------------------------
$ws.Range("A26:K32").Copy()
$ws.Range("A70").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 2. Fill in the new group's content.
# ---------------------------------------------------------------
# Group header
$ws.Range("A70").Value = "4."
$ws.Range("B70").Value = "Đăng ký tài khoản"

# Test case row (row 72) - function name / item / expected result
$ws.Range("A72").Value = 11
$ws.Range("B72").Value = "DangKy_01"
$ws.Range("D72").Value = "Nhập thông tin không đầy đủ"
$ws.Range("F72").Value = "thông báo thông tin điền không đầy đủ ,yêu cầu điền đầy đủ thông tin yêu cầu vào các mục"

# Sub test-case id (row 74) and trailing number (row 76)
$ws.Range("A74").Value = "4-1"
$ws.Range("A76").Value = 11

# ---------------------------------------------------------------
# 3. Update the summary formulas on rows 1-4 (ranges shrink because
#    11 fewer template rows now lie below the sheet's real data).
# ---------------------------------------------------------------
$ws.Range("I1").Formula = '=COUNTIF(H1:H767,"OK")'
$ws.Range("I2").Formula = '=COUNTIF(H2:H768,"Not OK")'
$ws.Range("I3").Formula = '=COUNTIF(H2:H768,"Untested")'
$ws.Range("I4").Formula = '=COUNTIF(H3:H769,"Result")'

# ---------------------------------------------------------------
# 4. Extend the print area to cover the new rows.
# ---------------------------------------------------------------
$ws.PageSetup.PrintArea = 'A1:K112'

# ---------------------------------------------------------------
# 5. Update the sheet view so the new rows are in focus/selected,
#    matching where the author was working.
# ---------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Range("F72:G76").Select()

$wb.Save()
